# Update odds figures on Sheet1 (Jogos_da_Semana_FlashScore_2024-10-12)
# per the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("N5").Value = 10
$ws.Range("N11").Value = 6.6
$ws.Range("H12").Value = 3.75
$ws.Range("J12").Value = 2.27
$ws.Range("K12").Value = 2.25
$ws.Range("L12").Value = 4.25
$ws.Range("Q12").Value = 1.62
$ws.Range("S12").Value = 1.31
$ws.Range("T12").Value = 3.26
$ws.Range("AA12").Value = 13
$ws.Range("AD12").Value = 7.5
$ws.Range("AH12").Value = 25
$ws.Range("AM12").Value = 350
$ws.Range("AP12").Value = 15.5
$ws.Range("AQ12").Value = 27
$ws.Range("AR12").Value = 50
$ws.Range("AT12").Value = 3.05
$ws.Range("AU12").Value = 6.9
$ws.Range("AW12").Value = 6
$ws.Range("AY12").Value = 25
$ws.Range("BB12").Value = 300
$ws.Range("G13").Value = 2.55
$ws.Range("H13").Value = 3.5
$ws.Range("I13").Value = 2.6
$ws.Range("J13").Value = 3
$ws.Range("L13").Value = 3.1
$ws.Range("Q13").Value = 1.6
$ws.Range("R13").Value = 2.3
$ws.Range("S13").Value = 1.3
$ws.Range("T13").Value = 3.4
$ws.Range("W13").Value = 12
$ws.Range("X13").Value = 15
$ws.Range("Y13").Value = 10
$ws.Range("Z13").Value = 26
$ws.Range("AA13").Value = 19
$ws.Range("AB13").Value = 23
$ws.Range("AC13").Value = 15
$ws.Range("AG13").Value = 12
$ws.Range("AH13").Value = 15
$ws.Range("AI13").Value = 10
$ws.Range("AJ13").Value = 26
$ws.Range("AK13").Value = 19
$ws.Range("AO13").Value = 13
$ws.Range("AP13").Value = 19
$ws.Range("AT13").Value = 3.4
$ws.Range("AX13").Value = 13
$ws.Range("AY13").Value = 19
$ws.Range("O14").Value = 1.57
$ws.Range("P14").Value = 2.25
$ws.Range("G23").Value = 2.5
$ws.Range("I23").Value = 2.7
$ws.Range("J23").Value = 3.25
$ws.Range("L23").Value = 3.4
$ws.Range("N23").Value = 10
$ws.Range("O23").Value = 1.36
$ws.Range("P23").Value = 3
$ws.Range("Q23").Value = 2.1
$ws.Range("R23").Value = 1.7
$ws.Range("AI23").Value = 10
$ws.Range("AJ23").Value = 26
$ws.Range("AQ23").Value = 51
$ws.Range("AX23").Value = 15
$ws.Range("BA23").Value = 67
$ws.Range("G24").Value = 2.05
$ws.Range("I24").Value = 4
$ws.Range("K24").Value = 2.1
$ws.Range("X24").Value = 9
$ws.Range("Z24").Value = 17
$ws.Range("AA24").Value = 17
$ws.Range("AO24").Value = 11
$ws.Range("M25").Value = 1.05
$ws.Range("N25").Value = 11
$ws.Range("Q25").Value = 1.98
$ws.Range("R25").Value = 1.88
$ws.Range("H31").Value = 3.6
$ws.Range("I31").Value = 4.35
$ws.Range("J31").Value = 2.27
$ws.Range("K31").Value = 2.18
$ws.Range("W31").Value = 7.1
$ws.Range("AA31").Value = 13.5
$ws.Range("AB31").Value = 26
$ws.Range("AG31").Value = 13
$ws.Range("AH31").Value = 26
$ws.Range("AJ31").Value = 75
$ws.Range("AL31").Value = 40
$ws.Range("AM31").Value = 500
$ws.Range("AO31").Value = 8.5
$ws.Range("AP31").Value = 18
$ws.Range("AU31").Value = 7.5
$ws.Range("AY31").Value = 30
$ws.Range("N35").Value = 1.03
$ws.Range("M37").Value = 1.03
$ws.Range("N37").Value = 7
